$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 19: Completed points becomes a static value instead of a formula ---
$ws.Range("D19").Value = 1.5

# --- Block 1: Sprint 2023.04.0 Bristol (rows 21-29) ---
$ws.Range("A21:D21").Merge()
$ws.Range("A21").Value = "Sprint 2023.04.0 Bristol"

$ws.Range("A22").Value = "Issue#"
$ws.Range("B22").Value = "Description"
$ws.Range("C22").Value = "Points"
$ws.Range("D22").Value = "Days worked:"

$ws.Range("A23").Value = 2864
$ws.Range("B23").Value = "Allow filtering on Evidence Type"
$ws.Range("C23").Value = 2
$ws.Range("D23").Value = 9

$ws.Range("A24").Value = 3223
$ws.Range("B24").Value = "Rename ""Delete"" to ""Move to trash bin"""
$ws.Range("C24").Value = 0.5
$ws.Range("D24").Value = "Availability:"

$ws.Range("A25").Value = 3227
$ws.Range("B25").Value = "Add Goals to My Evidence Table"
$ws.Range("C25").Value = 2
$ws.Range("D25").Value = 4.5

$ws.Range("D26").Value = "Planned points:"
$ws.Range("D27").Formula = "=SUM(C23:C25)"
$ws.Range("D28").Value = "Completed points:"
$ws.Range("D29").Formula = "=SUM(C23:C25)"

# --- Block 2: Sprint 2023.04.0 C**** (rows 31-39) ---
$ws.Range("A31:D31").Merge()
$ws.Range("A31").Value = "Sprint 2023.04.0 C****"

$ws.Range("A32").Value = "Issue#"
$ws.Range("B32").Value = "Description"
$ws.Range("C32").Value = "Points"
$ws.Range("D32").Value = "Days worked:"

$ws.Range("A33").Value = 2862
$ws.Range("D33").Value = 8

$ws.Range("A34").Value = 2861
$ws.Range("D34").Value = "Availability:"

$ws.Range("A35").Value = 2864
$ws.Range("D35").Value = 4

$ws.Range("D36").Value = "Planned points:"
$ws.Range("D37").Formula = "=SUM(C33:C35)"
$ws.Range("D38").Value = "Completed points:"
$ws.Range("D39").Formula = "=SUM(C33:C35)"

# --- Block 3: Sprint 2023.04.0 D**** (rows 41-49) ---
$ws.Range("A41:D41").Merge()
$ws.Range("A41").Value = "Sprint 2023.04.0 D****"

$ws.Range("A42").Value = "Issue#"
$ws.Range("B42").Value = "Description"
$ws.Range("C42").Value = "Points"
$ws.Range("D42").Value = "Days worked:"

$ws.Range("A43").Value = 2862
$ws.Range("D43").Value = 9

$ws.Range("A44").Value = 2861
$ws.Range("D44").Value = "Availability:"

$ws.Range("A45").Value = 2864
$ws.Range("D45").Value = 4

$ws.Range("D46").Value = "Planned points:"
$ws.Range("D47").Formula = "=SUM(C43:C45)"
$ws.Range("D48").Value = "Completed points:"
$ws.Range("D49").Formula = "=SUM(C43:C45)"

# --- Now apply the cell-level formatting on top of the values that are already in place ---
# (copying formats in after the content avoids disturbing formula evaluation of the new ranges)
$ws.Range("A11:D19").Copy()
$ws.Range("A21").PasteSpecial(-4122)

$ws.Range("A11:D19").Copy()
$ws.Range("A31").PasteSpecial(-4122)

$ws.Range("A11:D19").Copy()
$ws.Range("A41").PasteSpecial(-4122)

# --- Fix row 15 styling (becomes the "last row" of the first sprint's data, explanatory-text style) ---
$ws.Range("A8:C8").Copy()
$ws.Range("A15:C15").PasteSpecial(-4122)

# --- View state: scroll down, select B28 ---
$ws.Application.ActiveWindow.ScrollRow = 10
$ws.Range("B28").Select()
